$wb = $excel.ActiveWorkbook

# --- Rushing sheet updates (Week 16 log) ---
$rushing = $wb.Worksheets.Item("Rushing")

# C.McCoy
$rushing.Range("C2").Value = 31
$rushing.Range("D2").Value = 24

# C.Edmonds
$rushing.Range("C5").Value = 68
$rushing.Range("D5").Value = 30
$rushing.Range("E5").Value = 15
$rushing.Range("F5").Value = 18

# J.Ward
$rushing.Range("E7").Value = 4

# --- Receiving sheet updates (Week 16 log) ---
$receiving = $wb.Worksheets.Item("Receiving")

# C.Edmonds
$receiving.Range("C2").Value = 29
$receiving.Range("D2").Value = 23
$receiving.Range("G2").Value = 1
$receiving.Range("H2").Value = 1

# C.Kirk
$receiving.Range("C5").Value = 56
$receiving.Range("E5").Value = 22
$receiving.Range("F5").Value = 15

# Z.Ertz
$receiving.Range("C6").Value = 58
$receiving.Range("D6").Value = 47
$receiving.Range("E6").Value = 24

# A.Wesley
$receiving.Range("C8").Value = 14
$receiving.Range("D8").Value = 8
$receiving.Range("E8").Value = 5
$receiving.Range("F8").Value = 2
$receiving.Range("G8").Value = 3

# D.Daniels
$receiving.Range("C11").Value = 5
$receiving.Range("D11").Value = 3

# D.Harris
$receiving.Range("C12").Value = 62
$receiving.Range("D12").Value = 49
$receiving.Range("E12").Value = 10
$receiving.Range("G12").Value = 9

# --- Active sheet switches from Receiving back to Rushing ---
$rushing.Activate()
$rushing.Select()
$rushing.Range("A1").Select()
